$d = $word.ActiveDocument

# 1. "; prototypical examples include the Lisp and Scheme programming languages."
#    -> "; prototypical language examples include Lisp and Scheme."
$d.Content.Find.Execute(
    "; prototypical examples include the Lisp and Scheme programming languages.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "; prototypical language examples include Lisp and Scheme.", 2)

# 2. "...is a reader that correctly distinguishes" -> "...is a reader scheme that correctly distinguishes"
$d.Content.Find.Execute(
    "is a macro system for JavaScript whose primary contribution is a reader that correctly distinguishes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "is a macro system for JavaScript whose primary contribution is a reader scheme that correctly distinguishes", 2)

# 3. "...lexer and parser to eliminate the need..." -> "...lexer and parser and eliminates the need..."
$d.Content.Find.Execute(
    "This reader sits between the lexer and parser to eliminate the need for the bidirectional",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This reader sits between the lexer and parser and eliminates the need for the bidirectional", 2)

# 4. "...require additional structure to enable the use of macros." -> "...require additional structure around the token tree to enable the use of macros."
$d.Content.Find.Execute(
    "languages like JavaScript require additional structure to enable the use of macros.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "languages like JavaScript require additional structure around the token tree to enable the use of macros.", 2)

# 5. " (which were introduced by the Honu programming language)." -> " (which too were introduced in Honu)."
#    (kept outside the italic "custom operators" run so formatting is preserved)
$d.Content.Find.Execute(
    " (which were introduced by the Honu programming language).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (which too were introduced in Honu).", 2)

# 6. "by enabling matching of syntax before and after the macro identifier." -> "by enabling the matching of syntax both before and after the macro identifier."
$d.Content.Find.Execute(
    "by enabling matching of syntax before and after the macro identifier.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by enabling the matching of syntax both before and after the macro identifier.", 2)

# 7. "Sweet.js supports two primary styles namely " -> "Sweet.js supports two primary styles, namely "
#    (kept outside the italic "rule" run so formatting is preserved)
$d.Content.Find.Execute(
    "Sweet.js supports two primary styles namely ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sweet.js supports two primary styles, namely ", 2)

# 8. "in the paper is to have a macro read a file and to store the file contents into a string." -> "in the paper is one that reads a file and stores the file contents into a string."
$d.Content.Find.Execute(
    "in the paper is to have a macro read a file and to store the file contents into a string.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in the paper is one that reads a file and stores the file contents into a string.", 2)

# 9. Second "Disney et. al." mention: merge the italic ". " + "al." runs into
#    a single ". al." run (mirrors the already-merged first mention), and fix
#    the stray double space before "paper".  The search is scoped to start
#    after "While Disney " so only the second (paragraph 10) mention is hit,
#    and the italic "et" run is left untouched so its formatting survives.
$scopeAnchor = $d.Content
$foundAnchor = $scopeAnchor.Find.Execute("While Disney ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scopeStart = $scopeAnchor.End

$rMerge = $d.Range($scopeStart, $d.Content.End)
$rMerge.Find.Execute(". al.", $true, $false, $false, $false, $false, $true, 1, $false, ". al.", 2)

$rSpace = $d.Range($scopeStart, $d.Content.End)
$rSpace.Find.Execute("’s  paper only explicitly", $true, $false, $false, $false, $false, $true, 1, $false, "’s paper only explicitly", 2)

# 10. "can be applied to other ambiguous grammar languages." -> "can be applied to other languages with ambiguous grammars."
$d.Content.Find.Execute(
    "can be applied to other ambiguous grammar languages.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can be applied to other languages with ambiguous grammars.", 2)

# 11. Perl/Rust ambiguity sentence rewrite
$d.Content.Find.Execute(
    "Perl’s ambiguity around the forward slash (“/”) and Rust’s ambiguity when parsing the less than (“<”) symbol.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Perl’s similar ambiguity around forward slash (“/”) as well as Rust’s ambiguity when parsing the less than (“<”) symbol.", 2)
